$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new feature request rows (Crit/Crit-Fail sounds, Roll Again button)
$ws.Range("A20").Value = "Crit and Crit Fail Sounds"
$ws.Range("B20").Value = "When you roll minimum or maximum on a roll, do the willhelm scream or the air horn in order"
$ws.Range("D20").Value = "Stefan Titus"

$ws.Range("A21").Value = "Roll Again Button"
$ws.Range("B21").Value = "Stop the dismissing by random click for the results, make it so you have a roll again, and close buttons."
$ws.Range("D21").Value = "Weston Fiala"

# Re-apply the autofilter over the extended range (A1:E20), restoring the
# "blank" filter on the Completed Version column (colId 2 / Field 3). This
# also refreshes hidden rows for completed items, which fixes row 16
# (Horizontal custom rolls / Drop X High-Low area) that had been left visible
# despite having a Completed Version value.
$ws.AutoFilterMode = $false
$ws.Range("A1:E20").AutoFilter(3, @(""), 7)

# Keep the workbook-level _FilterDatabase defined name in sync with the
# autofilter's new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$20"
    }
}

# Update the selected/active cell to reflect where editing continued
$ws.Range("B24").Select()

$wb.Save()
